$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing text storage (prevents Excel from
# auto-converting numeric-looking strings like "322.34" into real numbers),
# then restore the cell style/format so no stray formatting is introduced.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "47.572.54"
$ws.Range("E2").Value = "  +4.55%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.490.95"
$ws.Range("E3").Value = "  +2.55%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.14%  "

# Row 5 - BNB
Set-TextValue "D5" "322.34"
$ws.Range("E5").Value = "  +1.14%  "

# Row 6 - Solana
Set-TextValue "D6" "104.99"
$ws.Range("E6").Value = "  +2.01%  "

# Row 7 - XRP
Set-TextValue "D7" "0.525"
$ws.Range("E7").Value = "  +1.60%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.07%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.542"
$ws.Range("E9").Value = "  +2.19%  "

# Row 10 - Avalanche
Set-TextValue "D10" "38.06"
$ws.Range("E10").Value = "  +6.81%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.0814"
$ws.Range("E11").Value = "  +1.07%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.23%  "

# Row 13 - Chainlink
Set-TextValue "D13" "18.31"
$ws.Range("E13").Value = "  +0.97%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.16"
$ws.Range("E14").Value = "  +1.37%  "

# Row 15 - Wrapped liquid staked Ether 2.0
Set-TextValue "D15" "2.882.48"
$ws.Range("E15").Value = "  +2.70%  "

# Row 16 - Wrapped Ether
Set-TextValue "D16" "2.494.42"
$ws.Range("E16").Value = "  +2.97%  "

# Row 17 - Polygon
Set-TextValue "D17" "0.848"
$ws.Range("E17").Value = "  +0.29%  "

# Row 18 - Wrapped BTC
Set-TextValue "D18" "47.456.17"
$ws.Range("E18").Value = "  +4.68%  "

# Row 19 - Internet Computer (DFINITY)
$ws.Range("E19").Value = "  +4.44%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +3.84%  "

# Row 21 - Shiba Inu
$ws.Range("E21").Value = "  +1.59%  "

# Row 22 - Litecoin
Set-TextValue "D22" "70.68"
$ws.Range("E22").Value = "  +2.56%  "

# Row 23 - Bitcoin Cash
Set-TextValue "D23" "251.54"
$ws.Range("E23").Value = "  +2.78%  "

# Row 24 - Immutable X
$ws.Range("E24").Value = "  +5.94%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  +3.10%  "

# Row 26 - Ethereum Classic
$ws.Range("E26").Value = "  +2.17%  "

# Row 27 - Dai
$ws.Range("E27").Value = "  +0.00%  "

# Row 28 - now Toncoin (was Cosmos)
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "2.28"
$ws.Range("E28").Value = "  +4.24%  "

# Row 29 - now Cosmos (was Toncoin)
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D29" "10.03"
$ws.Range("E29").Value = "  +4.44%  "

# Row 30 - Injective Protocol
$ws.Range("E30").Value = "  +6.45%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +7.89%  "

# Row 32 - OKB
Set-TextValue "D32" "49.42"
$ws.Range("E32").Value = "  +0.27%  "

# Row 33 - Celestia
Set-TextValue "D33" "19.76"
$ws.Range("E33").Value = "  -2.79%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  +2.89%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.20%  "

# Row 37 - ARBITRUM
Set-TextValue "D37" "1.97"
$ws.Range("E37").Value = "  +5.43%  "

# Row 38 - RenderToken
Set-TextValue "D38" "4.63"
$ws.Range("E38").Value = "  +4.21%  "

# Row 39 - LidoDAOToken
$ws.Range("E39").Value = "  +4.27%  "

# Row 40 - now Stellar (was WEMIXToken)
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D40" "0.112"
$ws.Range("E40").Value = "  +1.93%  "

# Row 41 - now WEMIXToken (was Stellar)
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D41" "2.25"
$ws.Range("E41").Value = "  +1.68%  "

# Row 42 - Monero
Set-TextValue "D42" "122.01"
$ws.Range("E42").Value = "  -3.12%  "

# Row 43 - EnergySwap
Set-TextValue "D43" "21.30"
$ws.Range("E43").Value = "  +3.96%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +2.20%  "

# Row 45 - Maker
Set-TextValue "D45" "1.967.26"
$ws.Range("E45").Value = "  +2.18%  "

# Row 46 - NEAR Protocol
$ws.Range("E46").Value = "  +1.60%  "

# Row 47 - ApeXProtocol
Set-TextValue "D47" "2.11"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48 - Stacks
$ws.Range("E48").Value = "  +0.81%  "

# Row 49 - FraxShare
Set-TextValue "D49" "9.21"
$ws.Range("E49").Value = "  +0.73%  "

# Row 50 - THORChain
$ws.Range("E50").Value = "  +11.60%  "

# Row 51 - BitcoinSV
Set-TextValue "D51" "79.52"
$ws.Range("E51").Value = "  +3.65%  "
